# Generate Report for Handback
# The file "6550b08e-f945-4eb6-81ea-8aec86d39a59.md" has been handed back and is now
# in sync with en-US. Update the Overview/status sheet and both per-locale handback
# report sheets (zh-cn, de-de) to reflect that: status text, target/handback file
# columns filled in, and the handback timestamp recorded.

$wb = $excel.ActiveWorkbook

$handedBackStatus = "Handed back: in sync with en-US"
$mdFile = "6550b08e-f945-4eb6-81ea-8aec86d39a59.md"

# --- Overview sheet: row for 6550b08e-....md now shows the handed-back status ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $handedBackStatus
$overview.Range("C2").Value = $handedBackStatus

function Update-LocaleSheet($SheetName, $HandoffXlf, $HandbackDatetime, $HandoffUrl, $HandbackUrl) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Status for the handed-back file.
    $ws.Range("B2").Value = $handedBackStatus

    # Latest Target File (E2) and Latest Handback File (F2) are now populated -
    # the target is the same markdown source, the handback file is the same
    # xlf that was handed off.
    $ws.Range("E2").Value = $mdFile
    $ws.Range("F2").Value = $HandoffXlf

    # Latest Handback DateTime (G2) is now recorded.
    $ws.Range("G2").Value = $HandbackDatetime

    # Style the two new link-like cells like the other file-name / file cells.
    $ws.Range("E2").Style = "HyperLink"
    $ws.Range("F2").Style = "HyperLink"

    # Add hyperlinks for the two newly-populated cells, matching the existing
    # hyperlink targets for the source file (A2) and handoff file (C2).
    $ws.Hyperlinks.Add($ws.Range("E2"), $HandoffUrl, [Type]::Missing, [Type]::Missing, $mdFile)
    $ws.Hyperlinks.Add($ws.Range("F2"), $HandbackUrl, [Type]::Missing, [Type]::Missing, $HandoffXlf)
}

Update-LocaleSheet "zh-cn" `
    "6550b08e-f945-4eb6-81ea-8aec86d39a59.b3da1c36f67cd3b60b347d584886aaa112f6ffa4.zh-cn.xlf" `
    "2016-03-09 08:08:56" `
    "https://github.com/OpenLocalizationTest/oltest/blob/5f876c8afd8280bcd41abdc918a3ca5e7562acdd/e2e/6550b08e-f945-4eb6-81ea-8aec86d39a59.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/825c9069dc05cc47a413b537f3c461bac503a3c9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6550b08e-f945-4eb6-81ea-8aec86d39a59.b3da1c36f67cd3b60b347d584886aaa112f6ffa4.zh-cn.xlf"

Update-LocaleSheet "de-de" `
    "6550b08e-f945-4eb6-81ea-8aec86d39a59.b3da1c36f67cd3b60b347d584886aaa112f6ffa4.de-de.xlf" `
    "2016-03-09 08:09:05" `
    "https://github.com/OpenLocalizationTest/oltest/blob/5f876c8afd8280bcd41abdc918a3ca5e7562acdd/e2e/6550b08e-f945-4eb6-81ea-8aec86d39a59.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/644d10915efe1ed065c710fc322fdafd77538217/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6550b08e-f945-4eb6-81ea-8aec86d39a59.b3da1c36f67cd3b60b347d584886aaa112f6ffa4.de-de.xlf"

Write-Host "Handback report generated."
